# QUBES Code First commit
# Updates the trailing "new data" rows on the TestData sheet:
#   - A4 now holds the next "NewData NN" placeholder value
#   - C5/C6 now hold the next Variant/WPL serial-number pair
# The previous run of this sheet (rows 4-6) carried the *last* entries from a
# much longer block of generated sample rows; this commit advances the
# generator to the next batch, so only the trailing cells actually change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 is being refreshed with newly generated data (it is no longer the
# bolded/styled "section" cell it used to be - the cell now carries plain,
# unformatted text), so strip its old formatting before writing the value.
$ws.Range("A4").ClearFormats()
$ws.Range("A4").Value = "NewData 69"

$ws.Range("C5").Value = "Var1-VS1P320220136"
$ws.Range("C6").Value = "WPL031081"
